$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" text cell A1 with new conversion rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.27 = 12467.32 pesos`n✅ 12467.32 pesos = 3.26 = 964.61 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update "tasas" sheet rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 306
$wsTasas.Range("O10").Value = 3815
$wsTasas.Range("N12").Value = 3827
$wsTasas.Range("O12").Value = 296.1
